$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024-Badge-BOM-REV1")

# Add comments for the fabricator: what not to place, and the 4-pin header.
$ws.Range("G12").Value = "DO NOT PLACE. No Solder paste"
$ws.Range("G15").Copy()
$ws.Range("G12").PasteSpecial(-4122)

$ws.Range("G17").Value = "4 pin-header"
$ws.Range("G15").Copy()
$ws.Range("G17").PasteSpecial(-4122)
